$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Status" column header to "Term" and refresh its values to the
# new payment-term taxonomy (Cash / 30 Days / 7 Days) in place of the old
# Active / Inactive / Pending status values.
$ws.Range("F1").Value = "Term"

# Write "Cash" before "30 Days" so the rebuilt shared-strings table orders
# the new unique strings the same way the source workbook does.
$ws.Range("F6").Value = "Cash"

$ws.Range("F2").Value = "30 Days"
$ws.Range("F3").Value = "30 Days"
$ws.Range("F4").Value = "30 Days"
$ws.Range("F5").Value = "30 Days"
$ws.Range("F7").Value = "30 Days"
$ws.Range("F8").Value = "30 Days"
$ws.Range("F9").Value = "7 Days"
$ws.Range("F10").Value = "30 Days"
$ws.Range("F11").Value = "30 Days"

# Move the active selection from E12 to I12.
$ws.Range("I12").Select() | Out-Null

# Match the printed page setup (Letter-ish A4/Letter paper code 9, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
